$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.667.48'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '1.598.99'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.29%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '211.84'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('E7').Value = '  +0.27%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.0618'
$c.ClearFormats()
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +0.41%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '19.57'
$c.ClearFormats()
$ws.Range('E10').Value = '  -0.89%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0837'
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '1.823.43'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.03'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.550.65'
$ws.Range('E14').Value = '  -2.46%  '
$ws.Range('E15').Value = '  +0.03%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '65.22'
$c.ClearFormats()
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '26.669.23'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('E19').Value = '  +0.34%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '209.13'
$c.ClearFormats()
$ws.Range('E20').Value = '  -0.46%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '7.05'
$c.ClearFormats()
$ws.Range('E21').Value = '  +4.81%  '
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('E23').Value = '  +0.54%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '8.97'
$c.ClearFormats()
$ws.Range('E24').Value = '  +0.49%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '145.28'
$c.ClearFormats()
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('E28').Value = '  -0.53%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0515'
$c.ClearFormats()
$ws.Range('E30').Value = '  +2.34%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  +0.26%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '2.95'
$c.ClearFormats()
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').Value = '1.279.65'
$ws.Range('E34').Value = '  -1.60%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.623'
$c.ClearFormats()
$ws.Range('E35').Value = '  -7.43%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.45'
$c.ClearFormats()
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('E38').Value = '  -1.13%  '
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  +19.05%  '
$ws.Range('E41').Value = '  +2.60%  '
$ws.Range('E42').Value = '  +0.38%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.784'
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.79%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '63.95'
$c.ClearFormats()
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').Value = '1.735.91'
$ws.Range('E45').Value = '  +0.00%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '91.32'
$c.ClearFormats()
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('E47').Value = '  -2.17%  '
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('E50').Value = '  -0.10%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '7.41'
$c.ClearFormats()
$ws.Range('E51').Value = '  -1.35%  '
